# Daily update at 8 AM UTC
# Appends the next day's row of data to the bottom of the sheet, and moves
# the "last row" date number format down to the newly added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the distinctive "last row" date format before it gets overwritten.
$lastRowFormat = $ws.Range("A28").NumberFormat

# The old last data row (28) is no longer the last row, so give it the same
# date format used by all the other (non-last) rows.
$ws.Range("A28").NumberFormat = $ws.Range("A27").NumberFormat

# Append the new row of data.
$ws.Range("A29").Value = 45769
$ws.Range("B29").Value = 116
$ws.Range("C29").Value = 116
$ws.Range("D29").Value = 120

# The new last row takes on the distinctive "last row" date format.
$ws.Range("A29").NumberFormat = $lastRowFormat
